$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "68.326.98"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +0.89%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.748.53"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -0.97%  "
$ws.Range("E4").Value = "  +0.03%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "595.52"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.26%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "166.66"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  -1.59%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.746.22"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -1.02%  "
$ws.Range("E8").Value = "  -0.06%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.519"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  -0.94%  "
$ws.Range("E10").Value = "  -3.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.48"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  -0.32%  "
$ws.Range("E12").Value = "  -1.10%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000258"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  -7.43%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.00"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -1.67%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.378.09"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -0.96%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.745.62"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -0.80%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "68.312.01"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +0.81%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "17.87"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -3.85%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.00"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -2.58%  "
$ws.Range("E20").Value = "  -0.25%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.68"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +1.20%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.18"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  -0.34%  "
$ws.Range("E23").Value = "  -3.11%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "84.41"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  +0.63%  "
$ws.Range("E25").Value = "  -3.07%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.18"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  -0.88%  "
$ws.Range("E27").Value = "  -1.87%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.09"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  -1.92%  "
$ws.Range("E29").Value = "  -0.06%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "3.895.66"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -0.97%  "
$ws.Range("E31").Value = "  -4.60%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.28"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  -4.43%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "29.83"
$ws.Range("D33").Style = "Normal"
$ws.Range("E33").Value = "  -2.44%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "2.17"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.44%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.19"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  +0.51%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "3.704.46"
$ws.Range("D37").Style = "Normal"
$ws.Range("E37").Value = "  -1.12%  "
$ws.Range("E38").Value = "  -2.80%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "3.37"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -11.27%  "
$ws.Range("E40").Value = "  +0.32%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "0.999"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "  -0.44%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "5.79"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  -0.42%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "1.00"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.02%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "8.58"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -1.09%  "
$ws.Range("E47").Value = "  -1.52%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "42.78"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  +9.27%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "45.93"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +0.15%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "146.50"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  +4.32%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "390.38"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -1.36%  "
